$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal (non-numeric-auto-detected) text value into a cell
# without leaving behind any formula residue or extra cell styles.
# Technique: write a formula that evaluates to the literal text, then
# copy/paste-special-values over itself so the formula collapses into a
# plain cached string value (t="s" shared-string cell), matching how a
# "typed in" text-like numeric string is stored.
function Set-TextValue {
    param($addr, $text)
    $r = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
}

# Column A (Listen_start)
Set-TextValue "A2" "4799.8950226"
Set-TextValue "A3" "4822.0594192"

# Column B (Video_play)
Set-TextValue "B2" "4800.025131"
Set-TextValue "B3" "4822.0235126"

# Column C (flash detection) - value contains an embedded CR+LF ("365" then a
# line break), same text used for both rows so they share one string entry.
$cr = [char]13
$lf = [char]10
$flashValue = "365$cr$lf"
Set-TextValue "C2" $flashValue
Set-TextValue "C3" $flashValue

# Column D (Video_pause)
Set-TextValue "D2" "4811.4082363"
Set-TextValue "D3" "4833.5408547"

# Column E (Listen_stop)
Set-TextValue "E2" "4813.4254384"
Set-TextValue "E3" "4836.8184474"

# Column F (start_diff) - plain numeric values
$ws.Range("F2").Value = 0.1301084000006085
$ws.Range("F3").Value = -0.03590660000008938

# The embedded line break in column C can make the engine auto-expand the
# row height on save; restore rows 2-3 to their normal auto height so no
# stray ht/customHeight attributes are introduced.
$ws.Rows("2:3").EntireRow.AutoFit() | Out-Null
